$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "1.00", "63.576.08").
# Force text storage so Excel does not coerce/round it, then restore the
# default (unstyled) cell style so no stray formatting is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.576.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.696.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.65%  "
$ws.Range("E9").Value = "  -3.73%  "
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.371"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -8.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.172.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.441.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("E16").Value = "  -3.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.693.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "345.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.92%  "
$ws.Range("E21").Value = "  -4.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.509"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.97%  "
$ws.Range("E24").Value = "  -2.31%  "
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.40%  "
$ws.Range("E28").Value = "  -4.99%  "
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.65"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.80%  "
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "345.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.948"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.624"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0567"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.60%  "
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "130.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.09%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0243"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.42%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0976"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.54%  "
